$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
# Row 2
$ws.Range("D2").Value = '''28.974.82'
$ws.Range("E2").Value = '  -1.13%  '

# Row 3
$ws.Range("D3").Value = '''1.823.26'
$ws.Range("E3").Value = '  -1.20%  '

# Row 4
$ws.Range("D4").Value = '''0.9992'
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''241.41'
$ws.Range("E5").Value = '  +0.17%  '

# Row 6
$ws.Range("D6").Value = '''0.6367'
$ws.Range("E6").Value = '  -5.35%  '

# Row 7
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = '''44.50'
$ws.Range("E8").Value = '  +5.59%  '

# Row 9
$ws.Range("D9").Value = '''0.07315'
$ws.Range("E9").Value = '  -1.81%  '

# Row 10
$ws.Range("D10").Value = '''0.2920'
$ws.Range("E10").Value = '  -0.99%  '

# Row 11
$ws.Range("D11").Value = '''22.75'
$ws.Range("E11").Value = '  -0.95%  '

# Row 12
$ws.Range("D12").Value = '''0.07657'
$ws.Range("E12").Value = '  -0.72%  '

# Row 13
$ws.Range("D13").Value = '''1.823.31'
$ws.Range("E13").Value = '  -1.06%  '

# Row 14
$ws.Range("D14").Value = '''4.970'
$ws.Range("E14").Value = '  -0.80%  '

# Row 15
$ws.Range("D15").Value = '''0.6628'
$ws.Range("E15").Value = '  -1.39%  '

# Row 16
$ws.Range("D16").Value = '''82.11'
$ws.Range("E16").Value = '  -4.51%  '

# Row 17
$ws.Range("D17").Value = '''6.041'
$ws.Range("E17").Value = '  -1.80%  '

# Row 18
$ws.Range("D18").Value = '''0.000008572'
$ws.Range("E18").Value = '  +2.98%  '

# Row 19
$ws.Range("D19").Value = '''28.871.90'
$ws.Range("E19").Value = '  -1.43%  '

# Row 20
$ws.Range("D20").Value = '''2.081.59'
$ws.Range("E20").Value = '  -0.36%  '

# Row 21
$ws.Range("D21").Value = '''12.36'
$ws.Range("E21").Value = '  -1.43%  '

# Row 22
$ws.Range("D22").Value = '''223.44'
$ws.Range("E22").Value = '  -2.64%  '

# Row 23
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("D24").Value = '''7.067'
$ws.Range("E24").Value = '  -1.70%  '

# Row 25
$ws.Range("E25").Value = '  +0.03%  '

# Row 26
$ws.Range("D26").Value = '''157.69'
$ws.Range("E26").Value = '  -2.13%  '

# Row 27
$ws.Range("D27").Value = '''8.423'
$ws.Range("E27").Value = '  -3.35%  '

# Row 28
$ws.Range("D28").Value = '''0.1370'
$ws.Range("E28").Value = '  -2.74%  '

# Row 29
$ws.Range("D29").Value = '''17.87'
$ws.Range("E29").Value = '  -0.89%  '

# Row 30
$ws.Range("D30").Value = '''1.504'
$ws.Range("E30").Value = '  -0.48%  '

# Row 31
$ws.Range("D31").Value = '''4.077'
$ws.Range("E31").Value = '  -2.15%  '

# Row 32
$ws.Range("D32").Value = '''1.199'
$ws.Range("E32").Value = '  +0.49%  '

# Row 33
$ws.Range("D33").Value = '''4.004'
$ws.Range("E33").Value = '  -1.69%  '

# Row 34
$ws.Range("D34").Value = '''0.05278'
$ws.Range("E34").Value = '  -0.62%  '

# Row 35
$ws.Range("D35").Value = '''0.7376'
$ws.Range("E35").Value = '  -2.70%  '

# Row 36
$ws.Range("E36").Value = '  -2.92%  '

# Row 37
$ws.Range("D37").Value = '''1.149'
$ws.Range("E37").Value = '  +0.91%  '

# Row 38
$ws.Range("D38").Value = '''2.645'
$ws.Range("E38").Value = '  -1.27%  '

# Row 39
$ws.Range("D39").Value = '''1.283.12'
$ws.Range("E39").Value = '  -2.95%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '''0.01780'
$ws.Range("E40").Value = '  -1.38%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.740'
$ws.Range("E41").Value = '  +0.50%  '

# Row 42
$ws.Range("D42").Value = '''6.335'
$ws.Range("E42").Value = '  +5.86%  '

# Row 43
$ws.Range("D43").Value = '''0.8941'
$ws.Range("E43").Value = '  -2.78%  '

# Row 44
$ws.Range("D44").Value = '''0.9999'
$ws.Range("E44").Value = '  -0.21%  '

# Row 45
$ws.Range("D45").Value = '''102.43'
$ws.Range("E45").Value = '  -0.72%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '''1.979.68'
$ws.Range("E46").Value = '  -0.48%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '''0.00000000122'
$ws.Range("E47").Value = '  -0.84%  '

# Row 48
$ws.Range("D48").Value = '''0.5134'
$ws.Range("E48").Value = '  -0.63%  '

# Row 49
$ws.Range("D49").Value = '''63.95'
$ws.Range("E49").Value = '  -0.13%  '

# Row 50
$ws.Range("D50").Value = '''1.724'
$ws.Range("E50").Value = '  -3.06%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.05821'
$ws.Range("E51").Value = '  -2.22%  '
